$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Judge-6 score (column O) for contestant #5 (row 6): 0 -> 10 ---
$ws.Range("O6").Value = 10

# --- Contestant #12 (row 13): fix judge-9 (col I) formula 9.5/2 -> 9.5 ---
$ws.Range("I13").Formula = "=9.5"
# Re-apply the original (non-highlighted) number format/style that cell
# carried before the highlight was added, copying it from a same-style donor
# cell (N3, style 15) so we don't mint a new style entry.
$ws.Range("N3").Copy()
$ws.Range("I13").PasteSpecial(-4122)   # xlPasteFormats

# --- Summary table (rows 22-36): column C "bonus" scores updated ---
$ws.Range("C22").Value = 9.5
$ws.Range("C23").Value = 10
$ws.Range("C24").Value = 9.5
$ws.Range("C25").Value = 9
$ws.Range("C26").Value = 10
$ws.Range("C27").Value = 9
$ws.Range("C29").Value = 10
$ws.Range("C30").Value = 9.5
$ws.Range("C31").Value = 9
$ws.Range("C32").Value = 10
$ws.Range("C33").Value = 8
$ws.Range("C34").Value = 8.5
$ws.Range("C35").Value = 9
$ws.Range("C36").Value = 9.5

# --- Active selection moved to O14 ---
[void]$ws.Range("O14").Select()
